$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 4).Value = '24.685.69'
$ws.Cells.Item(2, 5).Value = '  -0.04%  '

# Row 3
$ws.Cells.Item(3, 4).Value = '1.688.25'
$ws.Cells.Item(3, 5).Value = '  -0.95%  '

# Row 4
$ws.Cells.Item(4, 4).Value = '''1.002'
$ws.Cells.Item(4, 5).Value = '  +0.43%  '

# Row 5
$ws.Cells.Item(5, 4).Value = '''315.41'
$ws.Cells.Item(5, 5).Value = '  +0.37%  '

# Row 6
$ws.Cells.Item(6, 4).Value = '''1.003'
$ws.Cells.Item(6, 5).Value = '  +0.48%  '

# Row 7
$ws.Cells.Item(7, 4).Value = '''0.3939'
$ws.Cells.Item(7, 5).Value = '  -0.90%  '

# Row 8
$ws.Cells.Item(8, 4).Value = '''0.4053'
$ws.Cells.Item(8, 5).Value = '  -0.33%  '

# Row 9
$ws.Cells.Item(9, 2).Value = 'Polygon'
$ws.Cells.Item(9, 3).Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Cells.Item(9, 4).Value = '''1.486'
$ws.Cells.Item(9, 5).Value = '  -1.81%  '

# Row 10
$ws.Cells.Item(10, 2).Value = 'BinanceUSD'
$ws.Cells.Item(10, 3).Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
$ws.Cells.Item(10, 4).Value = '''1.002'
$ws.Cells.Item(10, 5).Value = '  +0.42%  '

# Row 11
$ws.Cells.Item(11, 4).Value = '''52.65'
$ws.Cells.Item(11, 5).Value = '  -1.26%  '

# Row 12
$ws.Cells.Item(12, 4).Value = '''0.08833'
$ws.Cells.Item(12, 5).Value = '  +0.17%  '

# Row 13
$ws.Cells.Item(13, 4).Value = '''7.240'
$ws.Cells.Item(13, 5).Value = '  -1.44%  '

# Row 14
$ws.Cells.Item(14, 4).Value = '''23.53'
$ws.Cells.Item(14, 5).Value = '  +0.35%  '

# Row 15
$ws.Cells.Item(15, 4).Value = '''8.014'
$ws.Cells.Item(15, 5).Value = '  +6.14%  '

# Row 16
$ws.Cells.Item(16, 4).Value = '''0.00001315'
$ws.Cells.Item(16, 5).Value = '  -0.82%  '

# Row 17
$ws.Cells.Item(17, 4).Value = '1.690.83'
$ws.Cells.Item(17, 5).Value = '  -0.58%  '

# Row 18
$ws.Cells.Item(18, 4).Value = '''99.57'
$ws.Cells.Item(18, 5).Value = '  -1.38%  '

# Row 19
$ws.Cells.Item(19, 4).Value = '''0.07008'
$ws.Cells.Item(19, 5).Value = '  -1.63%  '

# Row 20
$ws.Cells.Item(20, 4).Value = '''19.50'
$ws.Cells.Item(20, 5).Value = '  -0.35%  '

# Row 21
$ws.Cells.Item(21, 4).Value = '''6.992'
$ws.Cells.Item(21, 5).Value = '  +3.17%  '

# Row 22
$ws.Cells.Item(22, 4).Value = '''1.007'
$ws.Cells.Item(22, 5).Value = '  +0.94%  '

# Row 23
$ws.Cells.Item(23, 4).Value = '''14.28'
$ws.Cells.Item(23, 5).Value = '  +0.45%  '

# Row 24
$ws.Cells.Item(24, 4).Value = '24.668.12'
$ws.Cells.Item(24, 5).Value = '  -0.12%  '

# Row 25
$ws.Cells.Item(25, 4).Value = '''3.285'
$ws.Cells.Item(25, 5).Value = '  +9.22%  '

# Row 26
$ws.Cells.Item(26, 4).Value = '''2.356'
$ws.Cells.Item(26, 5).Value = '  +2.01%  '

# Row 27
$ws.Cells.Item(27, 4).Value = '''22.71'
$ws.Cells.Item(27, 5).Value = '  +1.13%  '

# Row 28
$ws.Cells.Item(28, 4).Value = '''162.68'
$ws.Cells.Item(28, 5).Value = '  +1.93%  '

# Row 29
$ws.Cells.Item(29, 2).Value = 'HuobiToken'
$ws.Cells.Item(29, 3).Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Cells.Item(29, 4).Value = '''5.187'
$ws.Cells.Item(29, 5).Value = '  +1.31%  '

# Row 30
$ws.Cells.Item(30, 2).Value = 'BitcoinCash'
$ws.Cells.Item(30, 3).Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Cells.Item(30, 4).Value = '''135.25'
$ws.Cells.Item(30, 5).Value = '  +1.12%  '

# Row 31
$ws.Cells.Item(31, 4).Value = '''7.613'
$ws.Cells.Item(31, 5).Value = '  +2.98%  '

# Row 32
$ws.Cells.Item(32, 4).Value = '1.875.09'
$ws.Cells.Item(32, 5).Value = '  -0.70%  '

# Row 33
$ws.Cells.Item(33, 4).Value = '''1.059'
$ws.Cells.Item(33, 5).Value = '  -2.61%  '

# Row 34
$ws.Cells.Item(34, 4).Value = '''0.08534'
$ws.Cells.Item(34, 5).Value = '  -1.94%  '

# Row 35
$ws.Cells.Item(35, 4).Value = '''7.103'
$ws.Cells.Item(35, 5).Value = '  -3.15%  '

# Row 36
$ws.Cells.Item(36, 4).Value = '''11.27'
$ws.Cells.Item(36, 5).Value = '  +2.71%  '

# Row 37
$ws.Cells.Item(37, 4).Value = '''0.2726'
$ws.Cells.Item(37, 5).Value = '  -0.25%  '

# Row 38
$ws.Cells.Item(38, 4).Value = '''1.887'
$ws.Cells.Item(38, 5).Value = '  -3.24%  '

# Row 39
$ws.Cells.Item(39, 4).Value = '''14.43'
$ws.Cells.Item(39, 5).Value = '  -2.89%  '

# Row 40
$ws.Cells.Item(40, 4).Value = '''0.09177'
$ws.Cells.Item(40, 5).Value = '  +1.78%  '

# Row 41
$ws.Cells.Item(41, 4).Value = '''0.02712'
$ws.Cells.Item(41, 5).Value = '  -2.58%  '

# Row 42
$ws.Cells.Item(42, 4).Value = '''1.464'
$ws.Cells.Item(42, 5).Value = '  -1.22%  '

# Row 43
$ws.Cells.Item(43, 4).Value = '''0.7622'
$ws.Cells.Item(43, 5).Value = '  -0.94%  '

# Row 44
$ws.Cells.Item(44, 4).Value = '''15.96'
$ws.Cells.Item(44, 5).Value = '  +3.10%  '

# Row 45
$ws.Cells.Item(45, 4).Value = '''2.586'
$ws.Cells.Item(45, 5).Value = '  +4.94%  '

# Row 46
$ws.Cells.Item(46, 4).Value = '''0.7134'
$ws.Cells.Item(46, 5).Value = '  -1.11%  '

# Row 47
$ws.Cells.Item(47, 4).Value = '''4.214'
$ws.Cells.Item(47, 5).Value = '  +0.93%  '

# Row 48
$ws.Cells.Item(48, 5).Value = '  +0.51%  '

# Row 49
$ws.Cells.Item(49, 2).Value = 'Flow'
$ws.Cells.Item(49, 3).Value = 'https://coinranking.com/coin/QQ0NCmjVq+flow-flow'
$ws.Cells.Item(49, 4).Value = '''1.318'
$ws.Cells.Item(49, 5).Value = '  +0.43%  '

# Row 50
$ws.Cells.Item(50, 2).Value = 'Quant'
$ws.Cells.Item(50, 3).Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Cells.Item(50, 4).Value = '''139.59'
$ws.Cells.Item(50, 5).Value = '  -1.66%  '

# Row 51
$ws.Cells.Item(51, 4).Value = '''0.07965'
$ws.Cells.Item(51, 5).Value = '  -0.69%  '
